$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.525.38"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.84"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.37"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.19"
$ws.Range("E8").Value = "  +5.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.05"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0889"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.794.08"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.557.92"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.505.72"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.19"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.23"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.09"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("E25").Value = "  +7.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.20"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.00"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.394.89"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -4.17%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.62"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.522"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0466"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.47"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.970"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.80"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.705.42"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.27"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  -0.82%  "
